# Auto-generated update of market-price / leve-profit columns (H:N) across 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2469
$ws.Range("I51").Value = 2316.6667
$ws.Range("J51").Value = 2926
$ws.Range("K51").Value = 2316.6667
$ws.Range("L51").Value = 2926
$ws.Range("M51").Value = -1832.6667
$ws.Range("N51").Value = -3894

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 19048550
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 21165010
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 63495030
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -63497246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5625.4062
$ws.Range("I116").Value = 6327.864
$ws.Range("J116").Value = 4080
$ws.Range("K116").Value = 6327.864
$ws.Range("L116").Value = 4080
$ws.Range("M116").Value = -2885.864
$ws.Range("N116").Value = -10964

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1799.6976
$ws.Range("I132").Value = 1596.5946
$ws.Range("J132").Value = 3052.1667
$ws.Range("K132").Value = 4789.783799999999
$ws.Range("L132").Value = 9156.500100000001
$ws.Range("M132").Value = -2259.783799999999
$ws.Range("N132").Value = -14216.5001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1573.48
$ws.Range("I137").Value = 1220.2972
$ws.Range("J137").Value = 2578.6924
$ws.Range("K137").Value = 3660.8916
$ws.Range("L137").Value = 7736.0772
$ws.Range("M137").Value = -1110.8916
$ws.Range("N137").Value = -12836.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 5204.273
$ws.Range("I41").Value = 2043.375
$ws.Range("J41").Value = 13633.333
$ws.Range("K41").Value = 2043.375
$ws.Range("L41").Value = 13633.333
$ws.Range("M41").Value = -1629.375
$ws.Range("N41").Value = -14461.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1722.4324
$ws.Range("I74").Value = 1403.6666
$ws.Range("K74").Value = 1403.6666
$ws.Range("M74").Value = -529.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1722.4324
$ws.Range("I77").Value = 1403.6666
$ws.Range("K77").Value = 7018.333000000001
$ws.Range("M77").Value = -2650.333000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 952239.75
$ws.Range("I122").Value = 1070936.4
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 3212809.2
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -3210359.2
$ws.Range("N122").Value = -12900.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2328503
$ws.Range("I132").Value = 1991.9722
$ws.Range("J132").Value = 14293416
$ws.Range("K132").Value = 5975.9166
$ws.Range("L132").Value = 42880248
$ws.Range("M132").Value = -3445.9166
$ws.Range("N132").Value = -42885308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9065.593999999999
$ws.Range("I20").Value = 1132.45
$ws.Range("J20").Value = 22287.5
$ws.Range("K20").Value = 1132.45
$ws.Range("L20").Value = 22287.5
$ws.Range("M20").Value = -885.45
$ws.Range("N20").Value = -22781.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1267.4138
$ws.Range("I94").Value = 760.7143
$ws.Range("J94").Value = 2597.5
$ws.Range("K94").Value = 760.7143
$ws.Range("L94").Value = 2597.5
$ws.Range("M94").Value = -309.7143
$ws.Range("N94").Value = -3499.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 167.5
$ws.Range("I22").Value = 125
$ws.Range("J22").Value = 181.66667
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = 181.66667
$ws.Range("M22").Value = 225
$ws.Range("N22").Value = -881.6666700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 220727.23
$ws.Range("I31").Value = 1650.4584
$ws.Range("J31").Value = 677930.9399999999
$ws.Range("K31").Value = 1650.4584
$ws.Range("L31").Value = 677930.9399999999
$ws.Range("M31").Value = -1355.4584
$ws.Range("N31").Value = -678520.9399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 220727.23
$ws.Range("I34").Value = 1650.4584
$ws.Range("J34").Value = 677930.9399999999
$ws.Range("K34").Value = 1650.4584
$ws.Range("L34").Value = 677930.9399999999
$ws.Range("M34").Value = -1448.4584
$ws.Range("N34").Value = -678334.9399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 4304.375
$ws.Range("I35").Value = 1995
$ws.Range("J35").Value = 8153.3335
$ws.Range("K35").Value = 1995
$ws.Range("L35").Value = 8153.3335
$ws.Range("M35").Value = -1701
$ws.Range("N35").Value = -8741.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2262.6667
$ws.Range("I122").Value = 2513.2727
$ws.Range("J122").Value = 1987
$ws.Range("K122").Value = 7539.8181
$ws.Range("L122").Value = 5961
$ws.Range("M122").Value = -5089.8181
$ws.Range("N122").Value = -10861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1777.4688
$ws.Range("I132").Value = 954.2
$ws.Range("J132").Value = 3149.5833
$ws.Range("K132").Value = 2862.6
$ws.Range("L132").Value = 9448.749899999999
$ws.Range("M132").Value = -332.6000000000004
$ws.Range("N132").Value = -14508.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4170934.2
$ws.Range("I33").Value = 7692403
$ws.Range("J33").Value = 9198.362999999999
$ws.Range("K33").Value = 46154418
$ws.Range("L33").Value = 55190.178
$ws.Range("M33").Value = -46154135
$ws.Range("N33").Value = -55756.178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 450.63635
$ws.Range("J107").Value = 485.125
$ws.Range("L107").Value = 1455.375
$ws.Range("N107").Value = -5295.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 13574.881
$ws.Range("I137").Value = 6809.6816
$ws.Range("J137").Value = 21016.6
$ws.Range("K137").Value = 20429.0448
$ws.Range("L137").Value = 63049.8
$ws.Range("M137").Value = -15329.0448
$ws.Range("N137").Value = -73249.79999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 266.4
$ws.Range("I2").Value = 230
$ws.Range("J2").Value = 284.6
$ws.Range("K2").Value = 230
$ws.Range("L2").Value = 284.6
$ws.Range("M2").Value = -117
$ws.Range("N2").Value = -510.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5504.08
$ws.Range("I70").Value = 5595.3184
$ws.Range("J70").Value = 4835
$ws.Range("K70").Value = 5595.3184
$ws.Range("L70").Value = 4835
$ws.Range("M70").Value = -5325.3184
$ws.Range("N70").Value = -5375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5504.08
$ws.Range("I73").Value = 5595.3184
$ws.Range("J73").Value = 4835
$ws.Range("K73").Value = 5595.3184
$ws.Range("L73").Value = 4835
$ws.Range("M73").Value = -4659.3184
$ws.Range("N73").Value = -6707

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 8298.117
$ws.Range("J136").Value = 8298.117
$ws.Range("L136").Value = 24894.351
$ws.Range("N136").Value = -29994.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1299.9231
$ws.Range("I16").Value = 1299.9231
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1299.9231
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1129.9231
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 40002708
$ws.Range("I40").Value = 55558376
$ws.Range("J40").Value = 2414.8572
$ws.Range("K40").Value = 55558376
$ws.Range("L40").Value = 2414.8572
$ws.Range("M40").Value = -55558240
$ws.Range("N40").Value = -2686.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1165.4286
$ws.Range("I61").Value = 1083.7
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 1083.7
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -881.7
$ws.Range("N61").Value = -3204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2478.6
$ws.Range("I68").Value = 2243.75
$ws.Range("J68").Value = 3418
$ws.Range("K68").Value = 2243.75
$ws.Range("L68").Value = 3418
$ws.Range("M68").Value = -1494.75
$ws.Range("N68").Value = -4916

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2478.6
$ws.Range("I71").Value = 2243.75
$ws.Range("J71").Value = 3418
$ws.Range("K71").Value = 11218.75
$ws.Range("L71").Value = 17090
$ws.Range("M71").Value = -7474.75
$ws.Range("N71").Value = -24578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 58015.85
$ws.Range("I82").Value = 1965.4546
$ws.Range("J82").Value = 126521.89
$ws.Range("K82").Value = 1965.4546
$ws.Range("L82").Value = 126521.89
$ws.Range("M82").Value = -1604.4546
$ws.Range("N82").Value = -127243.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 58015.85
$ws.Range("I85").Value = 1965.4546
$ws.Range("J85").Value = 126521.89
$ws.Range("K85").Value = 1965.4546
$ws.Range("L85").Value = 126521.89
$ws.Range("M85").Value = -717.4546
$ws.Range("N85").Value = -129017.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1165.4286
$ws.Range("I113").Value = 1083.7
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1083.7
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 1086.3
$ws.Range("N113").Value = -7140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2557567
$ws.Range("I122").Value = 3578673.8
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 10736021.4
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -10733571.4
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7580974
$ws.Range("I132").Value = 11116913
$ws.Range("J132").Value = 3962.0715
$ws.Range("K132").Value = 33350739
$ws.Range("L132").Value = 11886.2145
$ws.Range("M132").Value = -33348209
$ws.Range("N132").Value = -16946.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 23749.75
$ws.Range("I70").Value = 4999
$ws.Range("J70").Value = 30000
$ws.Range("K70").Value = 4999
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = -4684
$ws.Range("N70").Value = -30630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 23749.75
$ws.Range("I73").Value = 4999
$ws.Range("J73").Value = 30000
$ws.Range("K73").Value = 4999
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = -3907
$ws.Range("N73").Value = -32184
